$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.023.99'
$ws.Range('E2').Value = '  -1.41%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.517.82'
$ws.Range('E3').Value = '  -0.06%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '584.91'
$ws.Range('E5').Value = '  -1.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '132.14'
$ws.Range('E6').Value = '  -1.36%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.517.85'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -0.13%  '
$ws.Range('E9').Value = '  -1.34%  '
$ws.Range('E10').Value = '  +0.11%  '
$ws.Range('E11').Value = '  -0.26%  '
$ws.Range('E12').Value = '  -2.21%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.111.43'
$ws.Range('E13').Value = '  -0.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.58'
$ws.Range('E14').Value = '  -0.03%  '
$ws.Range('E15').Value = '  +1.27%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000179'
$ws.Range('E16').Value = '  -1.54%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.507.71'
$ws.Range('E17').Value = '  -0.40%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '64.039.02'
$ws.Range('E18').Value = '  -1.49%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.83'
$ws.Range('E19').Value = '  -2.77%  '
$ws.Range('E20').Value = '  -2.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.62'
$ws.Range('E21').Value = '  -1.50%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '384.17'
$ws.Range('E22').Value = '  -2.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.575'
$ws.Range('E23').Value = '  -0.72%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.658.10'
$ws.Range('E24').Value = '  -0.21%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '73.73'
$ws.Range('E25').Value = '  -1.34%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('E27').Value = '  +2.62%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.57'
$ws.Range('E28').Value = '  -0.52%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.52'
$ws.Range('E29').Value = '  -2.58%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('E31').Value = '  -0.69%  '
$ws.Range('E32').Value = '  -1.54%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.528.94'
$ws.Range('E33').Value = '  +0.09%  '
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.53'
$ws.Range('E35').Value = '  -2.29%  '
$ws.Range('E36').Value = '  +0.85%  '
$ws.Range('E37').Value = '  +0.56%  '
$ws.Range('E38').Value = '  -1.35%  '
$ws.Range('E39').Value = '  -1.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '159.06'
$ws.Range('E40').Value = '  -5.70%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0790'
$ws.Range('E41').Value = '  -2.91%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.814'
$ws.Range('E42').Value = '  -1.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '26.30'
$ws.Range('E43').Value = '  +1.37%  '
$ws.Range('E44').Value = '  -0.13%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '41.88'
$ws.Range('E45').Value = '  -2.27%  '
$ws.Range('E46').Value = '  -4.28%  '
$ws.Range('E47').Value = '  -0.57%  '
$ws.Range('E48').Value = '  -2.47%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.444.00'
$ws.Range('E49').Value = '  +1.21%  '
$ws.Range('E50').Value = '  -1.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.908'
$ws.Range('E51').Value = '  +0.16%  '
